# Estadisticos Matutinos 15 Oct
# Update columns E:K (Aprobados, Reprobados, %Aprobados, %Reprobados, Promedio,
# Aprobados(recup.), %Aprobados(recup.)) for the affected student rows in the
# "1er Parcial" and "3er Parcial" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 8
    $ws.Cells.Item(8, 5).Value = 20
    $ws.Cells.Item(8, 6).Value = 8
    $ws.Cells.Item(8, 7).Value = 71.43000000000001
    $ws.Cells.Item(8, 8).Value = 28.57
    $ws.Cells.Item(8, 9).Value = 7.1
    $ws.Cells.Item(8, 10).Value = 8
    $ws.Cells.Item(8, 11).Value = 28.57

    # Row 9
    $ws.Cells.Item(9, 5).Value = 21
    $ws.Cells.Item(9, 6).Value = 12
    $ws.Cells.Item(9, 7).Value = 63.64
    $ws.Cells.Item(9, 8).Value = 36.36
    $ws.Cells.Item(9, 9).Value = 6.8
    $ws.Cells.Item(9, 10).Value = 5
    $ws.Cells.Item(9, 11).Value = 15.15

    # Row 10
    $ws.Cells.Item(10, 5).Value = 24
    $ws.Cells.Item(10, 6).Value = 7
    $ws.Cells.Item(10, 7).Value = 77.42
    $ws.Cells.Item(10, 8).Value = 22.58
    $ws.Cells.Item(10, 9).Value = 7.7
    $ws.Cells.Item(10, 10).Value = 3
    $ws.Cells.Item(10, 11).Value = 9.68

    # Row 11
    $ws.Cells.Item(11, 5).Value = 10
    $ws.Cells.Item(11, 6).Value = 11
    $ws.Cells.Item(11, 7).Value = 47.62
    $ws.Cells.Item(11, 8).Value = 52.38
    $ws.Cells.Item(11, 9).Value = 7.2
    $ws.Cells.Item(11, 10).Value = 1
    $ws.Cells.Item(11, 11).Value = 4.76

    # Row 12
    $ws.Cells.Item(12, 5).Value = 29
    $ws.Cells.Item(12, 6).Value = 6
    $ws.Cells.Item(12, 7).Value = 82.86
    $ws.Cells.Item(12, 8).Value = 17.14
    $ws.Cells.Item(12, 9).Value = 8
    $ws.Cells.Item(12, 10).Value = 0
    $ws.Cells.Item(12, 11).Value = 0

    # Row 13
    $ws.Cells.Item(13, 5).Value = 12
    $ws.Cells.Item(13, 6).Value = 9
    $ws.Cells.Item(13, 7).Value = 57.14
    $ws.Cells.Item(13, 8).Value = 42.86
    $ws.Cells.Item(13, 9).Value = 6.9
    $ws.Cells.Item(13, 10).Value = 3
    $ws.Cells.Item(13, 11).Value = 14.29

    # Row 15
    $ws.Cells.Item(15, 5).Value = 25
    $ws.Cells.Item(15, 6).Value = 13
    $ws.Cells.Item(15, 7).Value = 65.79000000000001
    $ws.Cells.Item(15, 8).Value = 34.21
    $ws.Cells.Item(15, 9).Value = 9
    $ws.Cells.Item(15, 10).Value = 13
    $ws.Cells.Item(15, 11).Value = 34.21

    # Row 16
    $ws.Cells.Item(16, 5).Value = 28
    $ws.Cells.Item(16, 6).Value = 7
    $ws.Cells.Item(16, 7).Value = 80
    $ws.Cells.Item(16, 8).Value = 20
    $ws.Cells.Item(16, 9).Value = 9.199999999999999
    $ws.Cells.Item(16, 10).Value = 7
    $ws.Cells.Item(16, 11).Value = 20

    # Row 17
    $ws.Cells.Item(17, 5).Value = 21
    $ws.Cells.Item(17, 6).Value = 2
    $ws.Cells.Item(17, 7).Value = 91.3
    $ws.Cells.Item(17, 8).Value = 8.699999999999999
    $ws.Cells.Item(17, 9).Value = 9.5
    $ws.Cells.Item(17, 10).Value = 2
    $ws.Cells.Item(17, 11).Value = 8.699999999999999

    # Row 18
    $ws.Cells.Item(18, 5).Value = 16
    $ws.Cells.Item(18, 6).Value = 16
    $ws.Cells.Item(18, 7).Value = 50
    $ws.Cells.Item(18, 8).Value = 50
    $ws.Cells.Item(18, 9).Value = 8.800000000000001
    $ws.Cells.Item(18, 10).Value = 16
    $ws.Cells.Item(18, 11).Value = 50

    # Row 19
    $ws.Cells.Item(19, 5).Value = 27
    $ws.Cells.Item(19, 6).Value = 7
    $ws.Cells.Item(19, 7).Value = 79.41
    $ws.Cells.Item(19, 8).Value = 20.59
    $ws.Cells.Item(19, 9).Value = 9.699999999999999
    $ws.Cells.Item(19, 10).Value = 7
    $ws.Cells.Item(19, 11).Value = 20.59

    # Row 20
    $ws.Cells.Item(20, 5).Value = 22
    $ws.Cells.Item(20, 6).Value = 11
    $ws.Cells.Item(20, 7).Value = 66.67
    $ws.Cells.Item(20, 8).Value = 33.33
    $ws.Cells.Item(20, 9).Value = 9
    $ws.Cells.Item(20, 10).Value = 11
    $ws.Cells.Item(20, 11).Value = 33.33

    # Row 24
    $ws.Cells.Item(24, 5).Value = 22
    $ws.Cells.Item(24, 6).Value = 9
    $ws.Cells.Item(24, 7).Value = 70.97
    $ws.Cells.Item(24, 8).Value = 29.03
    $ws.Cells.Item(24, 9).Value = 6.6
    $ws.Cells.Item(24, 10).Value = 0
    $ws.Cells.Item(24, 11).Value = 0

    # Row 25
    $ws.Cells.Item(25, 5).Value = 29
    $ws.Cells.Item(25, 6).Value = 7
    $ws.Cells.Item(25, 7).Value = 80.56
    $ws.Cells.Item(25, 8).Value = 19.44
    $ws.Cells.Item(25, 9).Value = 7.5
    $ws.Cells.Item(25, 10).Value = 0
    $ws.Cells.Item(25, 11).Value = 0

    # Row 26
    $ws.Cells.Item(26, 5).Value = 34
    $ws.Cells.Item(26, 6).Value = 9
    $ws.Cells.Item(26, 7).Value = 79.06999999999999
    $ws.Cells.Item(26, 8).Value = 20.93
    $ws.Cells.Item(26, 9).Value = 7.5
    $ws.Cells.Item(26, 10).Value = 2
    $ws.Cells.Item(26, 11).Value = 4.65

    # Row 27
    $ws.Cells.Item(27, 5).Value = 38
    $ws.Cells.Item(27, 6).Value = 6
    $ws.Cells.Item(27, 7).Value = 86.36
    $ws.Cells.Item(27, 8).Value = 13.64
    $ws.Cells.Item(27, 9).Value = 7
    $ws.Cells.Item(27, 10).Value = 0
    $ws.Cells.Item(27, 11).Value = 0

    # Row 28
    $ws.Cells.Item(28, 5).Value = 21
    $ws.Cells.Item(28, 6).Value = 3
    $ws.Cells.Item(28, 7).Value = 87.5
    $ws.Cells.Item(28, 8).Value = 12.5
    $ws.Cells.Item(28, 9).Value = 7.1
    $ws.Cells.Item(28, 10).Value = 0
    $ws.Cells.Item(28, 11).Value = 0

    # Row 29
    $ws.Cells.Item(29, 5).Value = 23
    $ws.Cells.Item(29, 6).Value = 3
    $ws.Cells.Item(29, 7).Value = 88.45999999999999
    $ws.Cells.Item(29, 8).Value = 11.54
    $ws.Cells.Item(29, 9).Value = 6.7
    $ws.Cells.Item(29, 10).Value = 0
    $ws.Cells.Item(29, 11).Value = 0

    # Row 31
    $ws.Cells.Item(31, 5).Value = 31
    $ws.Cells.Item(31, 6).Value = 13
    $ws.Cells.Item(31, 7).Value = 70.45
    $ws.Cells.Item(31, 8).Value = 29.55
    $ws.Cells.Item(31, 9).Value = 7.1
    $ws.Cells.Item(31, 10).Value = 13
    $ws.Cells.Item(31, 11).Value = 29.55

    # Row 32
    $ws.Cells.Item(32, 5).Value = 30
    $ws.Cells.Item(32, 6).Value = 13
    $ws.Cells.Item(32, 7).Value = 69.77
    $ws.Cells.Item(32, 8).Value = 30.23
    $ws.Cells.Item(32, 9).Value = 7
    $ws.Cells.Item(32, 10).Value = 13
    $ws.Cells.Item(32, 11).Value = 30.23

    # Row 33
    $ws.Cells.Item(33, 5).Value = 12
    $ws.Cells.Item(33, 6).Value = 6
    $ws.Cells.Item(33, 7).Value = 66.67
    $ws.Cells.Item(33, 8).Value = 33.33
    $ws.Cells.Item(33, 9).Value = 6.9
    $ws.Cells.Item(33, 10).Value = 6
    $ws.Cells.Item(33, 11).Value = 33.33

    # Row 34
    $ws.Cells.Item(34, 5).Value = 15
    $ws.Cells.Item(34, 6).Value = 19
    $ws.Cells.Item(34, 7).Value = 44.12
    $ws.Cells.Item(34, 8).Value = 55.88
    $ws.Cells.Item(34, 9).Value = 9.5
    $ws.Cells.Item(34, 10).Value = 19
    $ws.Cells.Item(34, 11).Value = 55.88

    # Row 35
    $ws.Cells.Item(35, 5).Value = 15
    $ws.Cells.Item(35, 6).Value = 19
    $ws.Cells.Item(35, 7).Value = 44.12
    $ws.Cells.Item(35, 8).Value = 55.88
    $ws.Cells.Item(35, 9).Value = 9.5
    $ws.Cells.Item(35, 10).Value = 19
    $ws.Cells.Item(35, 11).Value = 55.88

}
